$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Gasolina"
$ws.Range("B7").Value = "Carro"
$ws.Range("C7").Value = "Gasto"
$ws.Range("D7").Value = 27000

$ws.Range("A8").Value = "Smile Direct Club"
$ws.Range("B8").Value = "Freelance"
$ws.Range("C8").Value = "Ingreso"
$ws.Range("D8").Value = 275000
